# Apply manually-fixed weight values (mean_Intake / sem_Intake columns P and Q)
# for several rows in the female early-session data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 (Session 1): refine existing mean/sem Intake values
$ws.Range("P19").Value = 59.82145489760461
$ws.Range("Q19").Value = 56.969700048919947

# Row 51 (Session 1): previously a sentinel 65535 with no sem value,
# now has real computed mean/sem Intake values.
$ws.Range("P51").Value = 56.220744202549618
$ws.Range("Q51").Value = 19.425478785356162

# Row 53 (Session 3): previously a sentinel 65535 with no sem value,
# now has real computed mean/sem Intake values.
$ws.Range("P53").Value = 114.90023813898675
$ws.Range("Q53").Value = 50.044468626942816

# Row 56 (Session 6): previously a sentinel 65535 with no sem value,
# now has real computed mean/sem Intake values.
$ws.Range("P56").Value = 61.087244108135195
$ws.Range("Q56").Value = 23.644092541978718

# Row 61 (Session 11): refine existing mean/sem Intake values
$ws.Range("P61").Value = 42.385652341041457
$ws.Range("Q61").Value = 9.7041622791667983

# Row 62 (Session 12): previously a sentinel 65535 with no sem value,
# now has real computed mean/sem Intake values.
$ws.Range("P62").Value = 37.897737798460241
$ws.Range("Q62").Value = 17.144928103836946
